$wb = $excel.ActiveWorkbook

# --- Rename sheet tabs (task order identifiers refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961144424143"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961160826461"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961160826461"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961161465964"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961162105973"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996114402416.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961144263744.csv"
$ws1.Range("B4").Value = "go_stims-16509961144263744.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961144424143.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509961159864101.csv"
$ws2.Range("B3").Value = "TB-1650996116066561.csv"
$ws2.Range("B4").Value = "OB-16509961153943777.csv"
$ws2.Range("B5").Value = "OB-1650996115186375.csv"
$ws2.Range("B6").Value = "ZB-match_7-16509961145143778.csv"
$ws2.Range("B7").Value = "OB-16509961147863748.csv"
$ws2.Range("B8").Value = "ZB-match_7-1650996114602421.csv"
$ws2.Range("B9").Value = "ZB-match_3-16509961144584117.csv"
$ws2.Range("B10").Value = "TB-16509961157784116.csv"

# --- Sheet 3: RS --- (no data changes, only tab name above)

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961161146102.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961160905957.csv"
$ws4.Range("B4").Value = "MM_stims-16509961161306055.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961161146102.csv"
$ws4.Range("B6").Value = "MM_stims-16509961161465964.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961161306055.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509961161945977.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961161625717.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961161465964.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961161786017.csv"
